# Seznam materialu pro Standu.
# Update the List1 "material bill" sheet:
#  - fix a typo in an item name ("ratky cerny" -> "kratky cerny")
#  - rename the "N plast" connector items to "N plast dutinka"
#  - insert a new "1 plast dutinka" line with a balance formula
#  - add a new "BASE162101A" line at the bottom
#  - add subtotal / computed columns (C, D, E) for rows 17-21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new "1 plast dutinka" row just above the "faston" row.
# This shifts the existing "faston" / "valcovy konektor" rows down by one
# (old row 21 -> 22, old row 22 -> 23).
$ws.Rows("21:21").Insert()

# New last row: BASE162101A (referenced first so it lands right after the
# existing strings in the shared-string table).
$ws.Cells.Item(24, 1).Value = "BASE162101A"
$ws.Cells.Item(24, 2).Value = 1

# Rename the plastic-connector rows to include "dutinka".
$ws.Cells.Item(18, 1).Value = "3 plast dutinka"
$ws.Cells.Item(19, 1).Value = "4 plast dutinka"
$ws.Cells.Item(20, 1).Value = "5 plast dutinka"

# New row: "1 plast dutinka", quantity computed from the balance formula below.
$ws.Cells.Item(21, 1).Value = "1 plast dutinka"

# Fix the "ratky cerny" typo -> "kratky cerny" (done last so the corrected
# string is appended at the very end of the shared-string table).
$ws.Cells.Item(8, 1).Value = "kratky cerny"

# Subtotal / doubled total for the short/long wire rows (8-17).
$ws.Range("C17").Formula = "=SUM(B8:B17)"
$ws.Range("D17").Formula = "=C17*2"

# Per-connector price/factor columns and computed totals.
$ws.Range("C18").Value = 2
$ws.Range("D18").Formula = "=C18*B18"

$ws.Range("C19").Value = 4
$ws.Range("D19").Formula = "=C19*B19"

$ws.Range("C20").Value = 4
$ws.Range("D20").Formula = "=C20*B20"

$ws.Range("E20").Formula = "=SUM(D18:D20)"

# Remaining quantity of "1 plast dutinka" needed to balance the totals.
$ws.Range("B21").Formula = "=D17-E20"

$ws.Range("D24").Select()
